$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# Pivot_2 rows got re-computed results: update expected_length (D) and coverage (E)
Set-TextValue "D6"  "0.0499266666666646"
Set-TextValue "E6"  "0.9468"

Set-TextValue "D9"  "0.0149779999999994"
Set-TextValue "E9"  "0.9508"

Set-TextValue "D12" "0.0014978000000001"
Set-TextValue "E12" "0.9476"

Set-TextValue "D15" "0.149779999999994"
Set-TextValue "E15" "0.9494"

Set-TextValue "D18" "0.0499266666666647"
Set-TextValue "E18" "0.9468"

Set-TextValue "D21" "0.0149779999999994"
Set-TextValue "E21" "0.9508"

Set-TextValue "D24" "0.0014978000000001"
Set-TextValue "E24" "0.9476"
